# Fix handling location problem - update Location column with corrected values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Calcutta Boys School Entally, Kolkata -> Mira Road East, Thane
$ws.Range("B4").Value = "Mira Road East, Thane"

# Row 5: Ahmedabad Palace Road, Bhopal -> Kandivali, Mumbai
$ws.Range("B5").Value = "Kandivali, Mumbai"

# Row 6: Vishakapatnam Central Suryabagh, Visakhapatnam -> Malad West, Mumbai
$ws.Range("B6").Value = "Malad West, Mumbai"

# Row 8: Muvattupuzha, Ernakulam -> Thane West, Thane
$ws.Range("B8").Value = "Thane West, Thane"

# Row 9: Kondapur, Hyderabad -> Mumbai Central, Mumbai
$ws.Range("B9").Value = "Mumbai Central, Mumbai"

# Update the active selection to B3
$ws.Range("B3").Select()
